# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45172 (2023-09-28) to 45175 (2023-10-01).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 288 }

$ws.Range("C2:C$lastRow").Value = 45175
